$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Drop the trailing three rows (old 25-27 "Requisitos" entries) ---
#     so the sheet shrinks from 27 to 24 rows, matching the new dimension A1:C24.
$ws.Rows.Item(27).Delete()
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(25).Delete()

# --- 2. Clear cells whose content was removed entirely ---
$ws.Cells.Item(14,2).Clear()
$ws.Cells.Item(14,3).Clear()
$ws.Cells.Item(16,2).Clear()
$ws.Cells.Item(16,3).Clear()
$ws.Cells.Item(22,2).Clear()
$ws.Cells.Item(22,3).Clear()
$ws.Cells.Item(23,1).Clear()
$ws.Cells.Item(24,1).Clear()

# --- 3. Update cells that already held a value (keeps their existing cell style) ---
$ws.Cells.Item(10,2).Value = '519033 - Carlos Yujiro Shigue'
$ws.Cells.Item(10,3).Value = '519033 - Carlos Yujiro Shigue'
$ws.Cells.Item(15,2).Value = '1176388 - Luiz Tadeu Fernandes Eleno'
$ws.Cells.Item(15,3).Value = '1176388 - Luiz Tadeu Fernandes Eleno'
$ws.Cells.Item(16,1).Value = 'Syllabus:'
$ws.Cells.Item(17,1).Value = 'Avaliação:'
$ws.Cells.Item(18,1).Value = 'Método:'
$ws.Cells.Item(18,2).Value = '7797767 - Viktor Pastoukhov'
$ws.Cells.Item(18,3).Value = '7797767 - Viktor Pastoukhov'
$ws.Cells.Item(19,1).Value = 'Critério:'
$ws.Cells.Item(20,1).Value = 'Norma de recuperação:'
$ws.Cells.Item(21,1).Value = 'Bibliografia:'
$ws.Cells.Item(21,2).Value = 'Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'
$ws.Cells.Item(21,3).Value = 'Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'
$ws.Cells.Item(22,1).Value = 'Requisitos:'
$ws.Cells.Item(23,2).Value = 'LOB1006 -  Cálculo IV  (Requisito)
'
$ws.Cells.Item(23,3).Value = 'LOB1006 -  Cálculo IV  (Requisito)
'
$ws.Cells.Item(24,2).Value = 'LOM3260 -  Computação Científica em Python  (Requisito)
'
$ws.Cells.Item(24,3).Value = 'LOM3260 -  Computação Científica em Python  (Requisito)
'

# --- 4. Create brand-new cells: copy the column formatting first, then set the value ---
$ws.Cells.Item(3,1).Copy()
$ws.Cells.Item(13,1).PasteSpecial(-4122)
$ws.Cells.Item(13,1).Value = 'Programa resumido:'
$ws.Cells.Item(3,1).Copy()
$ws.Cells.Item(14,1).PasteSpecial(-4122)
$ws.Cells.Item(14,1).Value = 'Short syllabus:'
$ws.Cells.Item(3,1).Copy()
$ws.Cells.Item(15,1).PasteSpecial(-4122)
$ws.Cells.Item(15,1).Value = 'Programa:'
$ws.Cells.Item(2,2).Copy()
$ws.Cells.Item(19,2).PasteSpecial(-4122)
$ws.Cells.Item(19,2).Value = 'Aulas expositivas, seminários e exercícios comentados.'
$ws.Cells.Item(2,3).Copy()
$ws.Cells.Item(19,3).PasteSpecial(-4122)
$ws.Cells.Item(19,3).Value = 'Aulas expositivas, seminários e exercícios comentados.'
$ws.Cells.Item(2,2).Copy()
$ws.Cells.Item(20,2).PasteSpecial(-4122)
$ws.Cells.Item(20,2).Value = 'Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2.'
$ws.Cells.Item(2,3).Copy()
$ws.Cells.Item(20,3).PasteSpecial(-4122)
$ws.Cells.Item(20,3).Value = 'Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2.'

# --- 5. Re-apply the custom row heights for the final layout ---
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(17).RowHeight = 15
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).RowHeight = 15
$ws.Rows.Item(23).RowHeight = 30
$ws.Rows.Item(24).RowHeight = 30

$excel.CutCopyMode = $false
